$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Neg_Change")
$ws2 = $wb.Worksheets.Item("Pos_Change")

# Sheet2 (Pos_Change) currently has data in rows 2-14; target only needs rows 2-13.
# Clear row 14 entirely since it is being removed.
$ws2.Rows.Item(14).ClearContents()

# --- Neg_Change (sheet1) data rows 2-13 ---
$ws1.Cells.Item(2, 1).Value = "COALINDIA"
$ws1.Cells.Item(2, 2).Value = 434
$ws1.Cells.Item(2, 3).Value = 435.8
$ws1.Cells.Item(2, 4).Value = 425.6
$ws1.Cells.Item(2, 5).Value = 429
$ws1.Cells.Item(2, 6).Value = 7104004
$ws1.Cells.Item(2, 7).Value = 14407896
$ws1.Cells.Item(2, 8).Value = -0.5069367519032619
$ws1.Cells.Item(2, 9).Value = "COALINDIA"
$ws1.Cells.Item(3, 1).Value = "SBILIFE"
$ws1.Cells.Item(3, 2).Value = 2090
$ws1.Cells.Item(3, 3).Value = 2102.1
$ws1.Cells.Item(3, 4).Value = 2071.7
$ws1.Cells.Item(3, 5).Value = 2076.9
$ws1.Cells.Item(3, 6).Value = 425540
$ws1.Cells.Item(3, 7).Value = 1051481
$ws1.Cells.Item(3, 8).Value = -0.5952946368027573
$ws1.Cells.Item(3, 9).Value = "SBILIFE"
$ws1.Cells.Item(4, 1).Value = "AMBUJACEM"
$ws1.Cells.Item(4, 2).Value = 540
$ws1.Cells.Item(4, 3).Value = 545.2
$ws1.Cells.Item(4, 4).Value = 533.55
$ws1.Cells.Item(4, 5).Value = 537
$ws1.Cells.Item(4, 6).Value = 1093928
$ws1.Cells.Item(4, 7).Value = 2484225
$ws1.Cells.Item(4, 8).Value = -0.5596501927160382
$ws1.Cells.Item(4, 9).Value = "AMBUJACEM"
$ws1.Cells.Item(5, 1).Value = "JUBLFOOD"
$ws1.Cells.Item(5, 2).Value = 526.5
$ws1.Cells.Item(5, 3).Value = 528.6
$ws1.Cells.Item(5, 4).Value = 520.5
$ws1.Cells.Item(5, 5).Value = 527.1
$ws1.Cells.Item(5, 6).Value = 784250
$ws1.Cells.Item(5, 7).Value = 1704055
$ws1.Cells.Item(5, 8).Value = -0.5397742443759151
$ws1.Cells.Item(5, 9).Value = "JUBLFOOD"
$ws1.Cells.Item(6, 1).Value = "LICHSGFIN"
$ws1.Cells.Item(6, 2).Value = 520.2
$ws1.Cells.Item(6, 3).Value = 526.8
$ws1.Cells.Item(6, 4).Value = 514.3
$ws1.Cells.Item(6, 5).Value = 521
$ws1.Cells.Item(6, 6).Value = 1358409
$ws1.Cells.Item(6, 7).Value = 3296969
$ws1.Cells.Item(6, 8).Value = -0.5879824772389428
$ws1.Cells.Item(6, 9).Value = "LICHSGFIN"
$ws1.Cells.Item(7, 1).Value = "PHOENIXLTD"
$ws1.Cells.Item(7, 2).Value = 1900
$ws1.Cells.Item(7, 3).Value = 1925.8
$ws1.Cells.Item(7, 4).Value = 1880
$ws1.Cells.Item(7, 5).Value = 1892.5
$ws1.Cells.Item(7, 6).Value = 303561
$ws1.Cells.Item(7, 7).Value = 671211
$ws1.Cells.Item(7, 8).Value = -0.5477413212834712
$ws1.Cells.Item(7, 9).Value = "PHOENIXLTD"
$ws1.Cells.Item(8, 1).Value = "CUMMINSIND"
$ws1.Cells.Item(8, 2).Value = 3970
$ws1.Cells.Item(8, 3).Value = 4008.5
$ws1.Cells.Item(8, 4).Value = 3929.1
$ws1.Cells.Item(8, 5).Value = 3970
$ws1.Cells.Item(8, 6).Value = 693509
$ws1.Cells.Item(8, 7).Value = 1559613
$ws1.Cells.Item(8, 8).Value = -0.5553326370067446
$ws1.Cells.Item(8, 9).Value = "CUMMINSIND"
$ws1.Cells.Item(9, 1).Value = "SUPREMEIND"
$ws1.Cells.Item(9, 2).Value = 3500.8
$ws1.Cells.Item(9, 3).Value = 3530.8
$ws1.Cells.Item(9, 4).Value = 3426.7
$ws1.Cells.Item(9, 5).Value = 3475
$ws1.Cells.Item(9, 6).Value = 205745
$ws1.Cells.Item(9, 7).Value = 502803
$ws1.Cells.Item(9, 8).Value = -0.5908039530392619
$ws1.Cells.Item(9, 9).Value = "SUPREMEIND"
$ws1.Cells.Item(10, 1).Value = "VOLTAS"
$ws1.Cells.Item(10, 2).Value = 1470
$ws1.Cells.Item(10, 3).Value = 1479.8
$ws1.Cells.Item(10, 4).Value = 1445.5
$ws1.Cells.Item(10, 5).Value = 1455
$ws1.Cells.Item(10, 6).Value = 403722
$ws1.Cells.Item(10, 7).Value = 802829
$ws1.Cells.Item(10, 8).Value = -0.4971257889289002
$ws1.Cells.Item(10, 9).Value = "VOLTAS"
$ws1.Cells.Item(11, 1).Value = "HUDCO"
$ws1.Cells.Item(11, 2).Value = 217.55
$ws1.Cells.Item(11, 3).Value = 219.47
$ws1.Cells.Item(11, 4).Value = 212.7
$ws1.Cells.Item(11, 5).Value = 215.3
$ws1.Cells.Item(11, 6).Value = 2302407
$ws1.Cells.Item(11, 7).Value = 4724520
$ws1.Cells.Item(11, 8).Value = -0.5126685885550278
$ws1.Cells.Item(11, 9).Value = "HUDCO"
$ws1.Cells.Item(12, 1).Value = "EXIDEIND"
$ws1.Cells.Item(12, 2).Value = 350.1
$ws1.Cells.Item(12, 3).Value = 351.8
$ws1.Cells.Item(12, 4).Value = 343
$ws1.Cells.Item(12, 5).Value = 345.85
$ws1.Cells.Item(12, 6).Value = 917928
$ws1.Cells.Item(12, 7).Value = 2008447
$ws1.Cells.Item(12, 8).Value = -0.54296628190836
$ws1.Cells.Item(12, 9).Value = "EXIDEIND"
$ws1.Cells.Item(13, 1).Value = "DELHIVERY"
$ws1.Cells.Item(13, 2).Value = 399.3
$ws1.Cells.Item(13, 3).Value = 399.6
$ws1.Cells.Item(13, 4).Value = 389.35
$ws1.Cells.Item(13, 5).Value = 395.1
$ws1.Cells.Item(13, 6).Value = 1735186
$ws1.Cells.Item(13, 7).Value = 3844603
$ws1.Cells.Item(13, 8).Value = -0.548669654578119
$ws1.Cells.Item(13, 9).Value = "DELHIVERY"

# --- Pos_Change (sheet2) data rows 2-13 ---
$ws2.Cells.Item(2, 1).Value = "SBIN"
$ws2.Cells.Item(2, 2).Value = 1022.8
$ws2.Cells.Item(2, 3).Value = 1029.5
$ws2.Cells.Item(2, 4).Value = 1016
$ws2.Cells.Item(2, 5).Value = 1028.5
$ws2.Cells.Item(2, 6).Value = 11016895
$ws2.Cells.Item(2, 7).Value = 7617922
$ws2.Cells.Item(2, 8).Value = 0.4461811239337972
$ws2.Cells.Item(2, 9).Value = "SBIN"
$ws2.Cells.Item(3, 1).Value = "TITAN"
$ws2.Cells.Item(3, 2).Value = 4231.6
$ws2.Cells.Item(3, 3).Value = 4267.6
$ws2.Cells.Item(3, 4).Value = 4210
$ws2.Cells.Item(3, 5).Value = 4233.5
$ws2.Cells.Item(3, 6).Value = 770707
$ws2.Cells.Item(3, 7).Value = 514767
$ws2.Cells.Item(3, 8).Value = 0.4971958186907863
$ws2.Cells.Item(3, 9).Value = "TITAN"
$ws2.Cells.Item(4, 1).Value = "HDFCBANK"
$ws2.Cells.Item(4, 2).Value = 947.7
$ws2.Cells.Item(4, 3).Value = 947.7
$ws2.Cells.Item(4, 4).Value = 932.6
$ws2.Cells.Item(4, 5).Value = 934.75
$ws2.Cells.Item(4, 6).Value = 32042827
$ws2.Cells.Item(4, 7).Value = 21071666
$ws2.Cells.Item(4, 8).Value = 0.5206594011123753
$ws2.Cells.Item(4, 9).Value = "HDFCBANK"
$ws2.Cells.Item(5, 1).Value = "HINDUNILVR"
$ws2.Cells.Item(5, 2).Value = 2415.1
$ws2.Cells.Item(5, 3).Value = 2422.2
$ws2.Cells.Item(5, 4).Value = 2365
$ws2.Cells.Item(5, 5).Value = 2389.5
$ws2.Cells.Item(5, 6).Value = 1519683
$ws2.Cells.Item(5, 7).Value = 1057665
$ws2.Cells.Item(5, 8).Value = 0.4368282962941952
$ws2.Cells.Item(5, 9).Value = "HINDUNILVR"
$ws2.Cells.Item(6, 1).Value = "BHARTIARTL"
$ws2.Cells.Item(6, 2).Value = 2047
$ws2.Cells.Item(6, 3).Value = 2049.9
$ws2.Cells.Item(6, 4).Value = 2014.9
$ws2.Cells.Item(6, 5).Value = 2025.3
$ws2.Cells.Item(6, 6).Value = 7388496
$ws2.Cells.Item(6, 7).Value = 4689135
$ws2.Cells.Item(6, 8).Value = 0.575662888784392
$ws2.Cells.Item(6, 9).Value = "BHARTIARTL"
$ws2.Cells.Item(7, 1).Value = "ITC"
$ws2.Cells.Item(7, 2).Value = 339
$ws2.Cells.Item(7, 3).Value = 339.5
$ws2.Cells.Item(7, 4).Value = 333.7
$ws2.Cells.Item(7, 5).Value = 334.3
$ws2.Cells.Item(7, 6).Value = 24922339
$ws2.Cells.Item(7, 7).Value = 16156143
$ws2.Cells.Item(7, 8).Value = 0.5425921273412844
$ws2.Cells.Item(7, 9).Value = "ITC"
$ws2.Cells.Item(8, 1).Value = "RELIANCE"
$ws2.Cells.Item(8, 2).Value = 1485
$ws2.Cells.Item(8, 3).Value = 1485.8
$ws2.Cells.Item(8, 4).Value = 1444.7
$ws2.Cells.Item(8, 5).Value = 1456.9
$ws2.Cells.Item(8, 6).Value = 13499760
$ws2.Cells.Item(8, 7).Value = 8883745
$ws2.Cells.Item(8, 8).Value = 0.5196023748993246
$ws2.Cells.Item(8, 9).Value = "RELIANCE"
$ws2.Cells.Item(9, 1).Value = "BRITANNIA"
$ws2.Cells.Item(9, 2).Value = 5930
$ws2.Cells.Item(9, 3).Value = 5970
$ws2.Cells.Item(9, 4).Value = 5896
$ws2.Cells.Item(9, 5).Value = 5907
$ws2.Cells.Item(9, 6).Value = 327794
$ws2.Cells.Item(9, 7).Value = 228799
$ws2.Cells.Item(9, 8).Value = 0.4326723455959161
$ws2.Cells.Item(9, 9).Value = "BRITANNIA"
$ws2.Cells.Item(10, 1).Value = "MPHASIS"
$ws2.Cells.Item(10, 2).Value = 2832
$ws2.Cells.Item(10, 3).Value = 2878.5
$ws2.Cells.Item(10, 4).Value = 2805
$ws2.Cells.Item(10, 5).Value = 2860
$ws2.Cells.Item(10, 6).Value = 314103
$ws2.Cells.Item(10, 7).Value = 219456
$ws2.Cells.Item(10, 8).Value = 0.4312800743657043
$ws2.Cells.Item(10, 9).Value = "MPHASIS"
$ws2.Cells.Item(11, 1).Value = "ALKEM"
$ws2.Cells.Item(11, 2).Value = 5907
$ws2.Cells.Item(11, 3).Value = 5933.5
$ws2.Cells.Item(11, 4).Value = 5813.5
$ws2.Cells.Item(11, 5).Value = 5865
$ws2.Cells.Item(11, 6).Value = 183102
$ws2.Cells.Item(11, 7).Value = 118161
$ws2.Cells.Item(11, 8).Value = 0.5495975829588443
$ws2.Cells.Item(11, 9).Value = "ALKEM"
$ws2.Cells.Item(12, 1).Value = "LUPIN"
$ws2.Cells.Item(12, 2).Value = 2190.4
$ws2.Cells.Item(12, 3).Value = 2205
$ws2.Cells.Item(12, 4).Value = 2157.3
$ws2.Cells.Item(12, 5).Value = 2177.4
$ws2.Cells.Item(12, 6).Value = 892652
$ws2.Cells.Item(12, 7).Value = 623831
$ws2.Cells.Item(12, 8).Value = 0.4309195920048859
$ws2.Cells.Item(12, 9).Value = "LUPIN"
$ws2.Cells.Item(13, 1).Value = "MCX"
$ws2.Cells.Item(13, 2).Value = 2249
$ws2.Cells.Item(13, 3).Value = 2303
$ws2.Cells.Item(13, 4).Value = 2239
$ws2.Cells.Item(13, 5).Value = 2293
$ws2.Cells.Item(13, 6).Value = 2810584
$ws2.Cells.Item(13, 7).Value = 1905600
$ws2.Cells.Item(13, 8).Value = 0.4749076406381192
$ws2.Cells.Item(13, 9).Value = "MCX"
